# Update the "Bugs" sheet:
#  - Row 15 ("GDI-D: ..." bug note) gets reworded to mention S2SB as well.
#  - Two new bug notes are appended as rows 16 and 17.
#  - Selection moves to the newly-added last row (A17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bugs")

$ws.Range("A15").Value = 'GDI-D *OR* S2SB: "Active Desktop Recovery" text screwed up'
$ws.Range("A16").Value = "S2SB: Is it set up same as GDI endianness wise?"
$ws.Range("A17").Value = "GDI-D: Rightmost horizontal line missing when size_in=0x0010"

$ws.Range("A17").Select()
